$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row-position swaps
# among HuobiToken/Aave and InjectiveProtocol/FTXToken) per the Nov 28 2023 refresh.

$ws.Range("D2").Value = "37.456.14"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "2.033.11"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.72"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.07"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0801"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "2.333.75"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.40"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.744"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "2.023.47"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "37.392.20"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.21"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.07"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.98"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.15"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  +5.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.77"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  +8.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.32"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.75"
$ws.Range("E37").Value = "  +8.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "1.475.05"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.83"
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.06"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.33"
$ws.Range("E45").Value = "  -6.00%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.22"
$ws.Range("E46").Value = "  +15.83%  "
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "2.221.86"
$ws.Range("E51").Value = "  +0.39%  "
